$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("G2").Value = 6265
$ws.Range("K2").Value = 2400
$ws.Range("K3").Value = 2308
$ws.Range("K4").Value = 488
$ws.Range("K5").Value = 151
$ws.Range("K6").Value = 2868
$ws.Range("G7").Value = 24708
$ws.Range("K7").Value = 8215

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 163
$ws.Range("K3").Value = 160
$ws.Range("K4").Value = 31
$ws.Range("K6").Value = 185
$ws.Range("K7").Value = 551

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 71
$ws.Range("K3").Value = 58
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K6").Value = 87
$ws.Range("K7").Value = 319

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 92
$ws.Range("K7").Value = 264

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 148

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 5
$ws.Range("K4").Value = 31
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 248
$ws.Range("K8").Value = 551
$ws.Range("K11").Value = 180
$ws.Range("K12").Value = 11
$ws.Range("K15").Value = 82
$ws.Range("K18").Value = 54
$ws.Range("K19").Value = 240
$ws.Range("K20").Value = 181
$ws.Range("K29").Value = 420
$ws.Range("K31").Value = 94
$ws.Range("K33").Value = 319
$ws.Range("K36").Value = 99
$ws.Range("K37").Value = 264
$ws.Range("K42").Value = 288
$ws.Range("K44").Value = 78
$ws.Range("K50").Value = 53
$ws.Range("K52").Value = 225
$ws.Range("K54").Value = 152
$ws.Range("G63").Value = 283
$ws.Range("K63").Value = 27
$ws.Range("K67").Value = 323
$ws.Range("K71").Value = 23
$ws.Range("K73").Value = 80
$ws.Range("K75").Value = 33
$ws.Range("K76").Value = 117
$ws.Range("K77").Value = 59
$ws.Range("K78").Value = 109
$ws.Range("K83").Value = 182
$ws.Range("K84").Value = 58
$ws.Range("K88").Value = 99
$ws.Range("K89").Value = 109
$ws.Range("K91").Value = 77
$ws.Range("K94").Value = 98
$ws.Range("K97").Value = 70
$ws.Range("K99").Value = 148
$ws.Range("G101").Value = 24708
$ws.Range("K101").Value = 8215

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 101
$ws.Range("K3").Value = 96
$ws.Range("K4").Value = 21
$ws.Range("K7").Value = 323

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 140
$ws.Range("K4").Value = 24
$ws.Range("K7").Value = 420

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K4").Value = 7
$ws.Range("K5").Value = 12
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 73
$ws.Range("K3").Value = 88
$ws.Range("K7").Value = 288

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 61
$ws.Range("K7").Value = 181

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 81
$ws.Range("K4").Value = 10
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 248

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 25
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("Andersonville")
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 5

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 63
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 11
